$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 48950.477
$ws.Range("J70").Value = 1350.5264
$ws.Range("L70").Value = 4051.5792
$ws.Range("N70").Value = -4591.5792

# Row 73
$ws.Range("H73").Value = 48950.477
$ws.Range("J73").Value = 1350.5264
$ws.Range("L73").Value = 4051.5792
$ws.Range("N73").Value = -5923.5792

# Row 129
$ws.Range("H129").Value = 941.2361
$ws.Range("I129").Value = 498.2857
$ws.Range("J129").Value = 1048.1552
$ws.Range("K129").Value = 1494.8571
$ws.Range("L129").Value = 3144.4656
$ws.Range("M129").Value = 3505.1429
$ws.Range("N129").Value = -13144.4656

# Row 132
$ws.Range("H132").Value = 2270.6135
$ws.Range("I132").Value = 1320.6072
$ws.Range("J132").Value = 3933.125
$ws.Range("K132").Value = 3961.8216
$ws.Range("L132").Value = 11799.375
$ws.Range("M132").Value = -1431.8216
$ws.Range("N132").Value = -16859.375

# Row 138
$ws.Range("H138").Value = 2944253
$ws.Range("I138").Value = 8697124
$ws.Range("J138").Value = 3896.8223
$ws.Range("K138").Value = 26091372
$ws.Range("L138").Value = 11690.4669
$ws.Range("M138").Value = -26086232
$ws.Range("N138").Value = -21970.4669


$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -40976

# Row 63
$ws.Range("H63").Value = 4365
$ws.Range("I63").Value = 4755
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 4755
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -4069
$ws.Range("N63").Value = -4372

# Row 66
$ws.Range("H66").Value = 4365
$ws.Range("I66").Value = 4755
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 23775
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -20343
$ws.Range("N66").Value = -21864

# Row 122
$ws.Range("H122").Value = 1271.2413
$ws.Range("I122").Value = 1138.8422
$ws.Range("J122").Value = 1522.8
$ws.Range("K122").Value = 3416.5266
$ws.Range("L122").Value = 4568.4
$ws.Range("M122").Value = -966.5266000000001
$ws.Range("N122").Value = -9468.4

# Row 123
$ws.Range("H123").Value = 28426
$ws.Range("J123").Value = 28426
$ws.Range("L123").Value = 28426
$ws.Range("N123").Value = -38226

# Row 131
$ws.Range("H131").Value = 29142.428
$ws.Range("J131").Value = 29142.428
$ws.Range("L131").Value = 29142.428
$ws.Range("N131").Value = -39222.428

# Row 132
$ws.Range("H132").Value = 1976.4
$ws.Range("I132").Value = 1320.6
$ws.Range("J132").Value = 4599.6
$ws.Range("K132").Value = 3961.8
$ws.Range("L132").Value = 13798.8
$ws.Range("M132").Value = -1431.8
$ws.Range("N132").Value = -18858.8


$ws = $wb.Worksheets.Item("BSM")
# Row 109
$ws.Range("H109").Value = 31097.285
$ws.Range("J109").Value = 31097.285
$ws.Range("L109").Value = 31097.285
$ws.Range("N109").Value = -33871.285

# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Row 112
$ws.Range("H112").Value = 26999.75
$ws.Range("J112").Value = 26999.75
$ws.Range("L112").Value = 26999.75
$ws.Range("N112").Value = -29953.75


$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1257.8334
$ws.Range("I58").Value = 1225.7826
$ws.Range("K58").Value = 1225.7826
$ws.Range("M58").Value = -1022.7826

# Row 86
$ws.Range("H86").Value = 4325
$ws.Range("I86").Value = 4357
$ws.Range("J86").Value = 4298.3335
$ws.Range("K86").Value = 4357
$ws.Range("L86").Value = 4298.3335
$ws.Range("M86").Value = -3234
$ws.Range("N86").Value = -6544.3335

# Row 89
$ws.Range("H89").Value = 4325
$ws.Range("I89").Value = 4357
$ws.Range("J89").Value = 4298.3335
$ws.Range("K89").Value = 21785
$ws.Range("L89").Value = 21491.6675
$ws.Range("M89").Value = -16169
$ws.Range("N89").Value = -32723.6675

# Row 134
$ws.Range("H134").Value = 1632.4584
$ws.Range("I134").Value = 1438.9048
$ws.Range("J134").Value = 2987.3333
$ws.Range("K134").Value = 4316.7144
$ws.Range("L134").Value = 8961.999899999999
$ws.Range("M134").Value = -1781.7144
$ws.Range("N134").Value = -14031.9999

# Row 136
$ws.Range("H136").Value = 1257.8334
$ws.Range("I136").Value = 1225.7826
$ws.Range("K136").Value = 3677.3478
$ws.Range("M136").Value = -1127.3478


$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 114100
$ws.Range("I15").Value = 5000
$ws.Range("K15").Value = 5000
$ws.Range("M15").Value = -4712

# Row 29
$ws.Range("H29").Value = 30003.5
$ws.Range("I29").Value = 30003.5
$ws.Range("K29").Value = 30003.5
$ws.Range("M29").Value = -29713.5

# Row 57
$ws.Range("H57").Value = 39936.75
$ws.Range("J57").Value = 39936.75
$ws.Range("L57").Value = 39936.75
$ws.Range("N57").Value = -41576.75

# Row 64
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

# Row 67
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

# Row 81
$ws.Range("H81").Value = 114100
$ws.Range("I81").Value = 5000
$ws.Range("K81").Value = 5000
$ws.Range("M81").Value = -4002

# Row 84
$ws.Range("H84").Value = 114100
$ws.Range("I84").Value = 5000
$ws.Range("K84").Value = 15000
$ws.Range("M84").Value = -10008

# Row 109
$ws.Range("H109").Value = 20284.5
$ws.Range("J109").Value = 20284.5
$ws.Range("L109").Value = 20284.5
$ws.Range("N109").Value = -22364.5

# Row 113
$ws.Range("H113").Value = 1101.875
$ws.Range("I113").Value = 764.3333
$ws.Range("J113").Value = 1746.2727
$ws.Range("K113").Value = 764.3333
$ws.Range("L113").Value = 1746.2727
$ws.Range("M113").Value = 1405.6667
$ws.Range("N113").Value = -6086.2727

# Row 122
$ws.Range("H122").Value = 3169.0588
$ws.Range("I122").Value = 2929.625
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 8788.875
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -6338.875
$ws.Range("N122").Value = -25900

# Row 136
$ws.Range("H136").Value = 13102.305
$ws.Range("J136").Value = 13102.305
$ws.Range("L136").Value = 39306.915
$ws.Range("N136").Value = -44406.915


$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 13894493
$ws.Range("I122").Value = 22731380
$ws.Range("J122").Value = 7956
$ws.Range("K122").Value = 68194140
$ws.Range("L122").Value = 23868
$ws.Range("M122").Value = -68191690
$ws.Range("N122").Value = -28768

# Row 132
$ws.Range("H132").Value = 5976.1714
$ws.Range("I132").Value = 5718
$ws.Range("J132").Value = 6722
$ws.Range("K132").Value = 17154
$ws.Range("L132").Value = 20166
$ws.Range("M132").Value = -14624
$ws.Range("N132").Value = -25226

# Row 136
$ws.Range("H136").Value = 3640.9375
$ws.Range("I136").Value = 2937.5
$ws.Range("J136").Value = 5751.25
$ws.Range("K136").Value = 8812.5
$ws.Range("L136").Value = 17253.75
$ws.Range("M136").Value = -6262.5
$ws.Range("N136").Value = -22353.75


$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 17860344
$ws.Range("I122").Value = 25002680
$ws.Range("J122").Value = 4505
$ws.Range("K122").Value = 75008040
$ws.Range("L122").Value = 13515
$ws.Range("M122").Value = -75005590
$ws.Range("N122").Value = -18415

# Row 123
$ws.Range("H123").Value = 21975.666
$ws.Range("J123").Value = 21975.666
$ws.Range("L123").Value = 21975.666
$ws.Range("N123").Value = -31775.666

# Row 125
$ws.Range("H125").Value = 51212.145
$ws.Range("J125").Value = 51212.145
$ws.Range("L125").Value = 51212.145
$ws.Range("N125").Value = -61052.145

# Row 132
$ws.Range("H132").Value = 1636.5714
$ws.Range("I132").Value = 1196.7142
$ws.Range("J132").Value = 2956.1428
$ws.Range("K132").Value = 3590.1426
$ws.Range("L132").Value = 8868.428400000001
$ws.Range("M132").Value = -1060.1426
$ws.Range("N132").Value = -13928.4284

